$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Annual -> Quarterly Bonus Target Percent (K), refreshed
# Last Bonus Allocation Percent (L), and recalculated Bonus Target
# Local Currency (M) / USD (N) amounts for every associate row (3-52).
$values = @{
    "K3" = 5
    "L3" = 110
    "M3" = 11000
    "K4" = 4.5
    "L4" = 90
    "M4" = 8550
    "K5" = 3.75
    "L5" = 100
    "M5" = 6750
    "K6" = 3
    "L6" = 85
    "M6" = 4500
    "K7" = 3
    "L7" = 115
    "M7" = 4350
    "K8" = 2.5
    "L8" = 115
    "M8" = 3000
    "K9" = 2.5
    "L9" = 100
    "M9" = 2875
    "K10" = 2.5
    "M10" = 2950
    "K11" = 2.5
    "L11" = 100
    "M11" = 2800
    "K12" = 3
    "M12" = 4560
    "K13" = 3.75
    "L13" = 90
    "M13" = 6562.5
    "K14" = 4.5
    "L14" = 95
    "M14" = 8325
    "K15" = 3
    "L15" = 105
    "M15" = 4650
    "K16" = 3
    "L16" = 110
    "M16" = 4470
    "K17" = 2.5
    "L17" = 100
    "M17" = 2750
    "K18" = 2.5
    "L18" = 95
    "M18" = 3125
    "K19" = 3
    "L19" = 110
    "M19" = 4410
    "K20" = 2.5
    "L20" = 110
    "M20" = 3050
    "K21" = 2.5
    "L21" = 105
    "M21" = 2975
    "K22" = 2.5
    "L22" = 110
    "M22" = 2900
    "K23" = 5
    "L23" = 85
    "M23" = 10750
    "K24" = 3.75
    "L24" = 90
    "M24" = 6825
    "K25" = 4.5
    "M25" = 8460
    "K26" = 3
    "L26" = 90
    "M26" = 4740
    "K27" = 3
    "L27" = 95
    "M27" = 4560
    "K28" = 2.5
    "L28" = 85
    "M28" = 3200
    "K29" = 2.5
    "L29" = 100
    "M29" = 3050
    "K30" = 2.5
    "L30" = 100
    "M30" = 2950
    "K31" = 3
    "L31" = 105
    "M31" = 4650
    "K32" = 2.5
    "L32" = 115
    "M32" = 3125
    "K33" = 3
    "M33" = 3150
    "N33" = 3987.33
    "K34" = 2.5
    "L34" = 110
    "M34" = 1950
    "N34" = 2468.35
    "K35" = 3.75
    "L35" = 90
    "M35" = 6937.5
    "K36" = 4.5
    "L36" = 95
    "M36" = 8640
    "K37" = 3
    "L37" = 100
    "M37" = 4650
    "K38" = 3
    "L38" = 115
    "M38" = 4560
    "K39" = 2.5
    "M39" = 3125
    "K40" = 2.5
    "L40" = 90
    "M40" = 3050
    "K41" = 3
    "L41" = 105
    "M41" = 4440
    "K42" = 2.5
    "M42" = 3250
    "K43" = 3.75
    "L43" = 85
    "M43" = 7050
    "K44" = 4.5
    "L44" = 95
    "M44" = 8775
    "K45" = 3
    "L45" = 115
    "M45" = 4800
    "K46" = 3
    "L46" = 115
    "M46" = 4680
    "K47" = 3
    "L47" = 110
    "M47" = 4590
    "K48" = 2.5
    "L48" = 90
    "M48" = 3200
    "K49" = 2.5
    "L49" = 105
    "M49" = 3100
    "K50" = 2.5
    "L50" = 110
    "M50" = 3000
    "K51" = 2.5
    "M51" = 2950
    "K52" = 2.5
    "M52" = 2875
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# These rows no longer carry a 'Last Bonus Allocation Percent' value.
$clearRefs = @("L10", "L33", "L42", "L51", "L52")
foreach ($ref in $clearRefs) {
    $ws.Range($ref).ClearContents()
}
